# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" (fund-holding detail, same shape as the
#    other quarterly sheets) right before the existing "2022-Q2" sheet.
# 2. Update the "总计" (summary) sheet: a new row for 2022-Q3 is inserted at
#    the top of the data (row 2), pushing the existing quarters down by one
#    row and re-numbering the running index in column A.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the new "2022-Q3" sheet before "2022-Q2"
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

# re-fetch a fresh reference to the template sheet (indices shifted after Add)
# "2021-Q3" has the same shape (1 header row + 3 data rows) as the new
# "2022-Q3" sheet, so its index column (A2:A4) has styling for all 3 rows.
$template = $wb.Worksheets.Item("2022-Q2")
$idxTemplate = $wb.Worksheets.Item("2021-Q3")

# Copy the header row (labels + style) from the template sheet so the new
# sheet's formatting (bold/centered header, borders) matches the others.
$template.Range("B1:H1").Copy($q3.Range("B1:H1"))
# Copy the style used on the index column (A2:A4) from a same-shaped sheet.
$idxTemplate.Range("A2:A4").Copy($q3.Range("A2:A4"))

# Columns B,D,E,F,G hold numeric-looking text (fund codes / percentages)
# that must stay text (e.g. "3.60" keeps its trailing zero) instead of being
# coerced to a number.
$q3.Range("B2:B4").NumberFormat = "@"
$q3.Range("D2:G4").NumberFormat = "@"

# -- data rows --
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "161724"
$q3.Range("C2").Value = "招商中证煤炭等权指数（LOF）A"
$q3.Range("D2").Value = "18.69"
$q3.Range("E2").Value = "94.42"
$q3.Range("F2").Value = "3.60"
$q3.Range("G2").Value = "0.6728"
$q3.Range("H2").Value = 4

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "013596"
$q3.Range("C3").Value = "招商中证煤炭等权指数（LOF）C"
$q3.Range("D3").Value = "1.38"
$q3.Range("E3").Value = "94.42"
$q3.Range("F3").Value = "3.60"
$q3.Range("G3").Value = "0.0497"
$q3.Range("H3").Value = 4

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "016347"
$q3.Range("C4").Value = "招商中证煤炭等权指数（LOF）E"
$q3.Range("D4").Value = "0.10"
$q3.Range("E4").Value = "94.42"
$q3.Range("F4").Value = "3.60"
$q3.Range("G4").Value = "0.0036"
$q3.Range("H4").Value = 4

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet with the new 2022-Q3 row,
# shifting the previous rows down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Push existing rows 2-5 down to rows 3-6, working from the bottom up so
# earlier writes are not overwritten, and carry the index-column (A) style
# down onto the newly used row 6.
$summary.Range("A5").Copy($summary.Range("A6"))
$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2020-Q4"
$summary.Range("C6").Value = 1
$summary.Range("D6").Value = 0.44

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q1"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 0.15

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q3"
$summary.Range("C4").Value = 3
$summary.Range("D4").Value = 0.07000000000000001

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 1.13

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.73

# ---------------------------------------------------------------------
# Restore the original active sheet ("2020-Q4" was the selected tab before
# the edit; adding a worksheet makes the new sheet active by default).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
